# edit.ps1 — apply "New crime data collected" update to the 7th Precinct
# CompStat weekly workbook (rows 15-30: Murder .. TOTAL crime-category table).
#
# The new week's report (12/19/2022 - 12/25/2022, Volume 29 Number 51)
# replaces last week's figures. Besides the header banner text, a block of
# cells in the summary table change value; a handful of them additionally
# flip between a numeric cell and the "no data" placeholder text cells
# (shared strings "0" / "***.*") used elsewhere in the sheet, so those get
# their formatting copied from an existing placeholder/numeric cell before
# the new value is written (Range.Copy carries the cell style/format; the
# value is then set explicitly so we don't depend on whatever the donor
# cell happens to hold).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Banner header text: volume/issue number and the week-covering dates ---
$ws.Range("A8").Value = "Volume 29   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/19/2022  Through  12/25/2022"

# --- Cells that change numeric/text "kind" (and therefore style) ---
$ws.Range("C14").Copy($ws.Range("F15"))
$ws.Range("F15").Value = "0"
$ws.Range("C14").Copy($ws.Range("D20"))
$ws.Range("D20").Value = "0"
$ws.Range("E14").Copy($ws.Range("E20"))
$ws.Range("E20").Value = "***.*"
$ws.Range("C14").Copy($ws.Range("C22"))
$ws.Range("C22").Value = "0"
$ws.Range("D16").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 1
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("D27").Value = "0"
$ws.Range("E14").Copy($ws.Range("E27"))
$ws.Range("E27").Value = "***.*"
$ws.Range("C14").Copy($ws.Range("D30"))
$ws.Range("D30").Value = "0"
$ws.Range("E14").Copy($ws.Range("E30"))
$ws.Range("E30").Value = "***.*"

# --- Remaining cells: same kind/style, value only changes ---
$updates = @(
    @{Cell="L15"; Value=21.428571428571},
    @{Cell="C16"; Value=1},
    @{Cell="E16"; Value=-66.666666666666},
    @{Cell="F16"; Value=7},
    @{Cell="G16"; Value=10},
    @{Cell="H16"; Value=-30},
    @{Cell="I16"; Value=175},
    @{Cell="J16"; Value=114},
    @{Cell="K16"; Value=53.508771929824},
    @{Cell="L16"; Value=54.867256637168},
    @{Cell="M16"; Value=34.615384615384},
    @{Cell="N16"; Value=-78.501228501228},
    @{Cell="C17"; Value=1},
    @{Cell="D17"; Value=3},
    @{Cell="E17"; Value=-66.666666666666},
    @{Cell="F17"; Value=13},
    @{Cell="G17"; Value=19},
    @{Cell="H17"; Value=-31.578947368421},
    @{Cell="I17"; Value=207},
    @{Cell="J17"; Value=202},
    @{Cell="K17"; Value=2.475247524752},
    @{Cell="L17"; Value=71.074380165289},
    @{Cell="M17"; Value=65.6},
    @{Cell="N17"; Value=-5.909090909090},
    @{Cell="C18"; Value=1},
    @{Cell="D18"; Value=3},
    @{Cell="E18"; Value=-66.666666666666},
    @{Cell="G18"; Value=20},
    @{Cell="H18"; Value=-55},
    @{Cell="I18"; Value=157},
    @{Cell="J18"; Value=125},
    @{Cell="K18"; Value=25.6},
    @{Cell="L18"; Value=-10.285714285714},
    @{Cell="M18"; Value=58.585858585858},
    @{Cell="N18"; Value=-59.640102827763},
    @{Cell="C19"; Value=6},
    @{Cell="D19"; Value=15},
    @{Cell="E19"; Value=-60},
    @{Cell="F19"; Value=48},
    @{Cell="G19"; Value=78},
    @{Cell="H19"; Value=-38.461538461538},
    @{Cell="I19"; Value=750},
    @{Cell="J19"; Value=586},
    @{Cell="K19"; Value=27.986348122866},
    @{Cell="L19"; Value=68.539325842696},
    @{Cell="M19"; Value=164.084507042254},
    @{Cell="N19"; Value=65.198237885462},
    @{Cell="M20"; Value=18.367346938775},
    @{Cell="N20"; Value=-85.089974293059},
    @{Cell="C21"; Value=9},
    @{Cell="D21"; Value=24},
    @{Cell="E21"; Value=-62.5},
    @{Cell="F21"; Value=79},
    @{Cell="G21"; Value=128},
    @{Cell="H21"; Value=-38.28125},
    @{Cell="I21"; Value=1366},
    @{Cell="J21"; Value=1094},
    @{Cell="K21"; Value=24.862888482632},
    @{Cell="L21"; Value=47.835497835497},
    @{Cell="M21"; Value=95.422031473533},
    @{Cell="N21"; Value=-40.505226480836},
    @{Cell="E22"; Value=-100},
    @{Cell="G22"; Value=5},
    @{Cell="H22"; Value=-60},
    @{Cell="J22"; Value=20},
    @{Cell="K22"; Value=30},
    @{Cell="M22"; Value=44.444444444444},
    @{Cell="C23"; Value=4},
    @{Cell="D23"; Value=4},
    @{Cell="E23"; Value=0},
    @{Cell="F23"; Value=12},
    @{Cell="H23"; Value=0},
    @{Cell="I23"; Value=172},
    @{Cell="J23"; Value=165},
    @{Cell="K23"; Value=4.242424242424},
    @{Cell="L23"; Value=-7.027027027027},
    @{Cell="M23"; Value=31.297709923664},
    @{Cell="C24"; Value=18},
    @{Cell="D24"; Value=46},
    @{Cell="E24"; Value=-60.869565217391},
    @{Cell="F24"; Value=112},
    @{Cell="G24"; Value=189},
    @{Cell="H24"; Value=-40.740740740740},
    @{Cell="I24"; Value=2162},
    @{Cell="J24"; Value=1637},
    @{Cell="K24"; Value=32.070861331704},
    @{Cell="L24"; Value=168.905472636816},
    @{Cell="M24"; Value=195.759233926129},
    @{Cell="C25"; Value=3},
    @{Cell="D25"; Value=5},
    @{Cell="E25"; Value=-40},
    @{Cell="F25"; Value=27},
    @{Cell="G25"; Value=24},
    @{Cell="H25"; Value=12.5},
    @{Cell="I25"; Value=442},
    @{Cell="J25"; Value=325},
    @{Cell="K25"; Value=36},
    @{Cell="L25"; Value=24.157303370786},
    @{Cell="M25"; Value=35.582822085889},
    @{Cell="F26"; Value=1},
    @{Cell="G26"; Value=3},
    @{Cell="H26"; Value=-66.666666666666},
    @{Cell="J26"; Value=17},
    @{Cell="K26"; Value=58.823529411764},
    @{Cell="L26"; Value=12.5},
    @{Cell="F27"; Value=3},
    @{Cell="G27"; Value=1},
    @{Cell="H27"; Value=200},
    @{Cell="I27"; Value=47},
    @{Cell="K27"; Value=4.444444444444},
    @{Cell="L27"; Value=74.074074074074},
    @{Cell="N28"; Value=-50},
    @{Cell="N29"; Value=-52.941176470588}
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
